$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4513
$ws1.Range("F3").Value = 2490
$ws1.Range("F6").Value = 59
$ws1.Range("F10").Value = 168
$ws1.Range("F11").Value = 170
$ws1.Range("F12").Value = 1683
$ws1.Range("F13").Value = 304
$ws1.Range("F14").Value = 3674
$ws1.Range("F15").Value = 17
$ws1.Range("F16").Value = 245

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4513
$ws4.Range("F3").Value = 2490
$ws4.Range("F7").Value = 59
$ws4.Range("F12").Value = 168
$ws4.Range("F13").Value = 170
$ws4.Range("F16").Value = 1683
$ws4.Range("F17").Value = 304
$ws4.Range("F18").Value = 3675
$ws4.Range("F19").Value = 17
$ws4.Range("F20").Value = 245
